$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$pairs = @(
    @("Before", "Trước"),
    @("During", "Trong khi xảy ra"),
    @("After", "Sau"),
    @("What's Happened?", "Điều Gì Đã Xảy Ra?"),
    @("What's the Worst?", "Điều Tồi Tệ Nhất?"),
    @("Cascadia Quake", "Động Đất Cascadia"),
    @("Tsunami Zone", "Khu Vực Sóng Thần"),
    @("If the dams failed", "Nếu đập bị vỡ")
)

$row = 94
foreach ($pair in $pairs) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}

$ws.Range("A94:B101").Style = "Normal"
$ws.Range("A94:B101").Select()
